# DFW_MSLibrary_TimeSheet - apply timesheet update for rows 27-28 (week of 6/27-7/03/21)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 27: update time spent and activity description ---
$ws.Range("B27").Value = 0.14583333333333334
$ws.Range("C27").Value = "UX Study Interviews, Sync with Kiran on UX study reporting and persona creation, work on three persona stories, email to MS Library SRIRs and Library Tools Support teams for onboarding Washington Post after going through provided information and providing access/content/usage information that is needed for the portal"
$ws.Rows.Item(27).RowHeight = 75

# --- Row 28: fill in previously-empty time spent and activity description ---
$ws.Range("B28").Value = 0.125
$ws.Range("C28").Value = "UX Study interviews & follow-up/debrief; Archives project work; Call with Omdia representative and continued research for Expert Insights playbook; Q1 Library/D&I quarterly sync meeting w/ Kiran with Global D&I employees to discuss reading lists/portal layout for ERG and D&I reading lists"
$ws.Rows.Item(28).RowHeight = 60

# E33 (=SUM(B27:B33)) and B40 (=SUM(B2:B38)) are formulas and recalculate automatically.

# --- Update the view/scroll position to match author's saved selection ---
$ws.Activate() | Out-Null
$ws.Range("B29").Select() | Out-Null
